$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the month label in C1
$ws.Range("C1").Value = "Август"

# The first worker row keeps its row, but now shows the second worker's name
$ws.Range("B4").Value = "Марков Евгений Викторович"

# All the other worker rows (5-12) are cleared out: numbers, names and the
# day-by-day shaded cells all go back to blank / unfilled, but keep their
# existing cell borders.
$range = $ws.Range("A5:AH12")
$range.ClearContents()
$range.Interior.ColorIndex = -4142
$range.Interior.Pattern = -4142
